# Fruta / hortaliza, semanal
# Adds one new weekly price record for "Macroferia Regional de Talca - Ciruela".
# The new record is inserted as row 71 (pushing the existing rows 71 and 72
# down to 72 and 73 respectively), matching how the source data feed keeps
# the most recent week at the top of this block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 71 - this shifts the old
# row 71 -> 72 and the old row 72 -> 73, carrying all of their values and
# formatting (incl. the date-formatted column D) along with them.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with this week's record.
$ws.Cells.Item(71, 1).Value = 5
$ws.Cells.Item(71, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(71, 3).Value = "Maule"
$ws.Cells.Item(71, 4).Value = 44585
$ws.Cells.Item(71, 5).Value = 7
$ws.Cells.Item(71, 6).Value = "Fruta"
$ws.Cells.Item(71, 7).Value = 100103
$ws.Cells.Item(71, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(71, 9).Value = 100103002
$ws.Cells.Item(71, 10).Value = "Ciruela"
$ws.Cells.Item(71, 11).Value = "Black Amber"
$ws.Cells.Item(71, 12).Value = "Primera"
$ws.Cells.Item(71, 13).Value = 300
$ws.Cells.Item(71, 14).Value = 9000
$ws.Cells.Item(71, 15).Value = 9000
$ws.Cells.Item(71, 16).Value = 9000
$ws.Cells.Item(71, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(71, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(71, 19).Value = 500
$ws.Cells.Item(71, 20).Value = 18
